$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = '2025-06-13T15:45:04+00:00'
$meta.Range("B15").Value = '4.0.1'

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")
$elem.Range("AJ10").Value = 'ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}'
$elem.Range("AJ11").Value = 'ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}'
$elem.Range("AJ2").Value = 'dom-2:If the resource is contained in another resource, it SHALL NOT contain nested Resources {contained.contained.empty()}
dom-3:If the resource is contained in another resource, it SHALL be referred to from elsewhere in the resource or SHALL refer to the containing resource {contained.where(((''#''+id in (%resource.descendants().reference | %resource.descendants().as(canonical) | %resource.descendants().as(uri) | %resource.descendants().as(url))) or descendants().where(reference = ''#'').exists() or descendants().where(as(canonical) = ''#'').exists() or descendants().where(as(canonical) = ''#'').exists()).not()).trace(''unmatched'', id).empty()}dom-4:If a resource is contained in another resource, it SHALL NOT have a meta.versionId or a meta.lastUpdated {contained.meta.versionId.empty() and contained.meta.lastUpdated.empty()}dom-5:If a resource is contained in another resource, it SHALL NOT have a security label {contained.meta.security.empty()}dom-6:A resource should have narrative for robust management {text.`div`.exists()}lst-1:A list can only have an emptyReason if it is empty {emptyReason.empty() or entry.empty()}lst-2:The deleted flag can only be used if the mode of the list is "changes" {mode = ''changes'' or entry.deleted.empty()}lst-3:An entry date can only be used if the mode of the list is "working" {mode = ''working'' or entry.date.empty()}'
$elem.Range("AJ24").Value = 'ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
'
$elem.Range("AJ8").Value = ''
$elem.Range("AL2").Value = ''
$elem.Range("O12").Value = 'Modifier extensions allow for extensions that *cannot* be safely ignored to be clearly distinguished from the vast majority of extensions which can be safely ignored.  This promotes interoperability by eliminating the need for implementers to prohibit the presence of extensions. For further information, see the [definition of modifier extensions](http://hl7.org/fhir/R4/extensibility.html#modifierExtension).'
$elem.Range("O27").Value = 'Modifier extensions allow for extensions that *cannot* be safely ignored to be clearly distinguished from the vast majority of extensions which can be safely ignored.  This promotes interoperability by eliminating the need for implementers to prohibit the presence of extensions. For further information, see the [definition of modifier extensions](http://hl7.org/fhir/R4/extensibility.html#modifierExtension).'
$elem.Range("Y6").Value = 'A human language.'
$elem.Range("Z14").Value = 'http://hl7.org/fhir/ValueSet/list-status|4.0.1'
$elem.Range("Z15").Value = 'http://hl7.org/fhir/ValueSet/list-mode|4.0.1'
